$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 5830
$ws.Cells.Item(2, 2).Value = 45828
$ws.Cells.Item(3, 1).Value = 5750
$ws.Cells.Item(3, 2).Value = 45828.01041666666
$ws.Cells.Item(4, 1).Value = 5700
$ws.Cells.Item(4, 2).Value = 45828.02083333334
$ws.Cells.Item(5, 1).Value = 5650
$ws.Cells.Item(5, 2).Value = 45828.03125
$ws.Cells.Item(6, 1).Value = 5630
$ws.Cells.Item(6, 2).Value = 45828.04166666666
$ws.Cells.Item(7, 1).Value = 5610
$ws.Cells.Item(7, 2).Value = 45828.05208333334
$ws.Cells.Item(8, 1).Value = 5580
$ws.Cells.Item(8, 2).Value = 45828.0625
$ws.Cells.Item(9, 1).Value = 5560
$ws.Cells.Item(9, 2).Value = 45828.07291666666
$ws.Cells.Item(10, 1).Value = 5520
$ws.Cells.Item(10, 2).Value = 45828.08333333334
$ws.Cells.Item(11, 1).Value = 5500
$ws.Cells.Item(11, 2).Value = 45828.09375
$ws.Cells.Item(12, 1).Value = 5500
$ws.Cells.Item(12, 2).Value = 45828.10416666666
$ws.Cells.Item(13, 1).Value = 5500
$ws.Cells.Item(13, 2).Value = 45828.11458333334
$ws.Cells.Item(14, 1).Value = 5500
$ws.Cells.Item(14, 2).Value = 45828.125
$ws.Cells.Item(15, 1).Value = 5500
$ws.Cells.Item(15, 2).Value = 45828.13541666666
$ws.Cells.Item(16, 1).Value = 5500
$ws.Cells.Item(16, 2).Value = 45828.14583333334
$ws.Cells.Item(17, 1).Value = 5500
$ws.Cells.Item(17, 2).Value = 45828.15625
$ws.Cells.Item(18, 1).Value = 5500
$ws.Cells.Item(18, 2).Value = 45828.16666666666
$ws.Cells.Item(19, 1).Value = 5500
$ws.Cells.Item(19, 2).Value = 45828.17708333334
$ws.Cells.Item(20, 1).Value = 5500
$ws.Cells.Item(20, 2).Value = 45828.1875
$ws.Cells.Item(21, 1).Value = 5520
$ws.Cells.Item(21, 2).Value = 45828.19791666666
$ws.Cells.Item(22, 1).Value = 5600
$ws.Cells.Item(22, 2).Value = 45828.20833333334
$ws.Cells.Item(23, 1).Value = 5660
$ws.Cells.Item(23, 2).Value = 45828.21875
$ws.Cells.Item(24, 1).Value = 5730
$ws.Cells.Item(24, 2).Value = 45828.22916666666
$ws.Cells.Item(25, 1).Value = 5810
$ws.Cells.Item(25, 2).Value = 45828.23958333334
$ws.Cells.Item(26, 1).Value = 5920
$ws.Cells.Item(26, 2).Value = 45828.25
$ws.Cells.Item(27, 1).Value = 5990
$ws.Cells.Item(27, 2).Value = 45828.26041666666
$ws.Cells.Item(28, 1).Value = 6060
$ws.Cells.Item(28, 2).Value = 45828.27083333334
$ws.Cells.Item(29, 1).Value = 6110
$ws.Cells.Item(29, 2).Value = 45828.28125
$ws.Cells.Item(30, 1).Value = 6110
$ws.Cells.Item(30, 2).Value = 45828.29166666666
$ws.Cells.Item(31, 1).Value = 6110
$ws.Cells.Item(31, 2).Value = 45828.30208333334
$ws.Cells.Item(32, 1).Value = 6110
$ws.Cells.Item(32, 2).Value = 45828.3125
$ws.Cells.Item(33, 1).Value = 6070
$ws.Cells.Item(33, 2).Value = 45828.32291666666
$ws.Cells.Item(34, 1).Value = 6010
$ws.Cells.Item(34, 2).Value = 45828.33333333334
$ws.Cells.Item(35, 1).Value = 5950
$ws.Cells.Item(35, 2).Value = 45828.34375
$ws.Cells.Item(36, 1).Value = 5860
$ws.Cells.Item(36, 2).Value = 45828.35416666666
$ws.Cells.Item(37, 1).Value = 5780
$ws.Cells.Item(37, 2).Value = 45828.36458333334
$ws.Cells.Item(38, 1).Value = 5690
$ws.Cells.Item(38, 2).Value = 45828.375
$ws.Cells.Item(39, 1).Value = 5610
$ws.Cells.Item(39, 2).Value = 45828.38541666666
$ws.Cells.Item(40, 1).Value = 5540
$ws.Cells.Item(40, 2).Value = 45828.39583333334
$ws.Cells.Item(41, 1).Value = 5480
$ws.Cells.Item(41, 2).Value = 45828.40625
$ws.Cells.Item(42, 1).Value = 5410
$ws.Cells.Item(42, 2).Value = 45828.41666666666
$ws.Cells.Item(43, 1).Value = 5370
$ws.Cells.Item(43, 2).Value = 45828.42708333334
$ws.Cells.Item(44, 1).Value = 5350
$ws.Cells.Item(44, 2).Value = 45828.4375
$ws.Cells.Item(45, 1).Value = 5330
$ws.Cells.Item(45, 2).Value = 45828.44791666666
$ws.Cells.Item(46, 1).Value = 5320
$ws.Cells.Item(46, 2).Value = 45828.45833333334
$ws.Cells.Item(47, 1).Value = 5320
$ws.Cells.Item(47, 2).Value = 45828.46875
$ws.Cells.Item(48, 1).Value = 5320
$ws.Cells.Item(48, 2).Value = 45828.47916666666
$ws.Cells.Item(49, 1).Value = 5330
$ws.Cells.Item(49, 2).Value = 45828.48958333334
$ws.Cells.Item(50, 1).Value = 5330
$ws.Cells.Item(50, 2).Value = 45828.5
$ws.Cells.Item(51, 1).Value = 5340
$ws.Cells.Item(51, 2).Value = 45828.51041666666
$ws.Cells.Item(52, 1).Value = 5340
$ws.Cells.Item(52, 2).Value = 45828.52083333334
$ws.Cells.Item(53, 1).Value = 5350
$ws.Cells.Item(53, 2).Value = 45828.53125
$ws.Cells.Item(54, 1).Value = 5350
$ws.Cells.Item(54, 2).Value = 45828.54166666666
$ws.Cells.Item(55, 1).Value = 5350
$ws.Cells.Item(55, 2).Value = 45828.55208333334
$ws.Cells.Item(56, 1).Value = 5370
$ws.Cells.Item(56, 2).Value = 45828.5625
$ws.Cells.Item(57, 1).Value = 5400
$ws.Cells.Item(57, 2).Value = 45828.57291666666
$ws.Cells.Item(58, 1).Value = 5430
$ws.Cells.Item(58, 2).Value = 45828.58333333334
$ws.Cells.Item(59, 1).Value = 5460
$ws.Cells.Item(59, 2).Value = 45828.59375
$ws.Cells.Item(60, 1).Value = 5490
$ws.Cells.Item(60, 2).Value = 45828.60416666666
$ws.Cells.Item(61, 1).Value = 5530
$ws.Cells.Item(61, 2).Value = 45828.61458333334
$ws.Cells.Item(62, 1).Value = 5590
$ws.Cells.Item(62, 2).Value = 45828.625
$ws.Cells.Item(63, 1).Value = 5660
$ws.Cells.Item(63, 2).Value = 45828.63541666666
$ws.Cells.Item(64, 1).Value = 5730
$ws.Cells.Item(64, 2).Value = 45828.64583333334
$ws.Cells.Item(65, 1).Value = 5810
$ws.Cells.Item(65, 2).Value = 45828.65625
$ws.Cells.Item(66, 1).Value = 5870
$ws.Cells.Item(66, 2).Value = 45828.66666666666
$ws.Cells.Item(67, 1).Value = 5960
$ws.Cells.Item(67, 2).Value = 45828.67708333334
$ws.Cells.Item(68, 1).Value = 6040
$ws.Cells.Item(68, 2).Value = 45828.6875
$ws.Cells.Item(69, 1).Value = 6130
$ws.Cells.Item(69, 2).Value = 45828.69791666666
$ws.Cells.Item(70, 1).Value = 6280
$ws.Cells.Item(70, 2).Value = 45828.70833333334
$ws.Cells.Item(71, 1).Value = 6390
$ws.Cells.Item(71, 2).Value = 45828.71875
$ws.Cells.Item(72, 1).Value = 6510
$ws.Cells.Item(72, 2).Value = 45828.72916666666
$ws.Cells.Item(73, 1).Value = 6620
$ws.Cells.Item(73, 2).Value = 45828.73958333334
$ws.Cells.Item(74, 1).Value = 6730
$ws.Cells.Item(74, 2).Value = 45828.75
$ws.Cells.Item(75, 1).Value = 6830
$ws.Cells.Item(75, 2).Value = 45828.76041666666
$ws.Cells.Item(76, 1).Value = 6900
$ws.Cells.Item(76, 2).Value = 45828.77083333334
$ws.Cells.Item(77, 1).Value = 6940
$ws.Cells.Item(77, 2).Value = 45828.78125
$ws.Cells.Item(78, 1).Value = 6970
$ws.Cells.Item(78, 2).Value = 45828.79166666666
$ws.Cells.Item(79, 1).Value = 7020
$ws.Cells.Item(79, 2).Value = 45828.80208333334
$ws.Cells.Item(80, 1).Value = 7070
$ws.Cells.Item(80, 2).Value = 45828.8125
$ws.Cells.Item(81, 1).Value = 7160
$ws.Cells.Item(81, 2).Value = 45828.82291666666
$ws.Cells.Item(82, 1).Value = 7220
$ws.Cells.Item(82, 2).Value = 45828.83333333334
$ws.Cells.Item(83, 1).Value = 7250
$ws.Cells.Item(83, 2).Value = 45828.84375
$ws.Cells.Item(84, 1).Value = 7240
$ws.Cells.Item(84, 2).Value = 45828.85416666666
$ws.Cells.Item(85, 1).Value = 7180
$ws.Cells.Item(85, 2).Value = 45828.86458333334
$ws.Cells.Item(86, 1).Value = 7050
$ws.Cells.Item(86, 2).Value = 45828.875
$ws.Cells.Item(87, 1).Value = 6950
$ws.Cells.Item(87, 2).Value = 45828.88541666666
$ws.Cells.Item(88, 1).Value = 6830
$ws.Cells.Item(88, 2).Value = 45828.89583333334
$ws.Cells.Item(89, 1).Value = 6660
$ws.Cells.Item(89, 2).Value = 45828.90625
$ws.Cells.Item(90, 1).Value = 6470
$ws.Cells.Item(90, 2).Value = 45828.91666666666
$ws.Cells.Item(91, 1).Value = 6310
$ws.Cells.Item(91, 2).Value = 45828.92708333334
$ws.Cells.Item(92, 1).Value = 6210
$ws.Cells.Item(92, 2).Value = 45828.9375
$ws.Cells.Item(93, 1).Value = 6100
$ws.Cells.Item(93, 2).Value = 45828.94791666666
$ws.Cells.Item(94, 1).Value = 5860
$ws.Cells.Item(94, 2).Value = 45828.95833333334
$ws.Cells.Item(95, 1).Value = 5820
$ws.Cells.Item(95, 2).Value = 45828.96875
$ws.Cells.Item(96, 1).Value = 5760
$ws.Cells.Item(96, 2).Value = 45828.97916666666
$ws.Cells.Item(97, 1).Value = 5660
$ws.Cells.Item(97, 2).Value = 45828.98958333334
